$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row total (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row (B12): 66 -> 110, and corresponding fraction label (E12)
$ws.Range("B12").Value = 110
$ws.Range("E12").Value = "110/140"
